$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New block: IOP Timers "MODE" register (rows 250-262), appended after the
# existing "IRQMASTER" block which ends at row 248 (row 249 stays blank as a
# separator, matching the existing layout convention used throughout the
# sheet).
# ---------------------------------------------------------------------------

# Write the new field-name values (column C) first, in the exact order they
# were first authored, so the shared-strings table grows in that order:
# SyncMode, SyncEnable, ResetMode, IrqOnTarget, IrqOnOF, IrqToggle,
# IrqRepeat, ClockSrc, IrqRequest, ReachTarget, ReachOF, Prescale0, Prescale1.
$ws.Range("C251").Value = "SyncMode"
$ws.Range("C250").Value = "SyncEnable"
$ws.Range("C252").Value = "ResetMode"
$ws.Range("C253").Value = "IrqOnTarget"
$ws.Range("C254").Value = "IrqOnOF"
$ws.Range("C256").Value = "IrqToggle"
$ws.Range("C255").Value = "IrqRepeat"
$ws.Range("C257").Value = "ClockSrc"
$ws.Range("C259").Value = "IrqRequest"
$ws.Range("C260").Value = "ReachTarget"
$ws.Range("C261").Value = "ReachOF"
$ws.Range("C258").Value = "Prescale0"
$ws.Range("C262").Value = "Prescale1"

# Column A - register field-name header for the new block (reuses the
# existing "MODE" shared string already used elsewhere in the sheet).
$ws.Range("A250").Value = "MODE"

# Column B - field index (0-based order within the register).
$ws.Range("B250").Value = 0
$ws.Range("B251").Value = 1
$ws.Range("B252").Value = 2
$ws.Range("B253").Value = 3
$ws.Range("B254").Value = 4
$ws.Range("B255").Value = 5
$ws.Range("B256").Value = 6
$ws.Range("B257").Value = 7
$ws.Range("B258").Value = 8
$ws.Range("B259").Value = 9
$ws.Range("B260").Value = 10
$ws.Range("B261").Value = 11
$ws.Range("B262").Value = 12

# Column D - bit start position.
$ws.Range("D250").Value = 0
$ws.Range("D251").Value = 1
$ws.Range("D252").Value = 3
$ws.Range("D253").Value = 4
$ws.Range("D254").Value = 5
$ws.Range("D255").Value = 6
$ws.Range("D256").Value = 7
$ws.Range("D257").Value = 8
$ws.Range("D258").Value = 9
$ws.Range("D259").Value = 10
$ws.Range("D260").Value = 11
$ws.Range("D261").Value = 12
$ws.Range("D262").Value = 13

# Column E - bit length.
$ws.Range("E250").Value = 1
$ws.Range("E251").Value = 2
$ws.Range("E252").Value = 1
$ws.Range("E253").Value = 1
$ws.Range("E254").Value = 1
$ws.Range("E255").Value = 1
$ws.Range("E256").Value = 1
$ws.Range("E257").Value = 1
$ws.Range("E258").Value = 1
$ws.Range("E259").Value = 1
$ws.Range("E260").Value = 1
$ws.Range("E261").Value = 1
$ws.Range("E262").Value = 2

# Column F - initial value (all zero for this block).
$ws.Range("F250").Value = 0
$ws.Range("F251").Value = 0
$ws.Range("F252").Value = 0
$ws.Range("F253").Value = 0
$ws.Range("F254").Value = 0
$ws.Range("F255").Value = 0
$ws.Range("F256").Value = 0
$ws.Range("F257").Value = 0
$ws.Range("F258").Value = 0
$ws.Range("F259").Value = 0
$ws.Range("F260").Value = 0
$ws.Range("F261").Value = 0
$ws.Range("F262").Value = 0

# Column H - register forward declaration (only present on the header row
# of the block, same convention as every previous register block).
$ws.Range("H250").Formula = '="class IOPTimersTimerRegister_"&A250&"_t;"'

# Column I - field-index constexpr declaration. Row 250 gets its own
# formula (first cell of the block); rows 251:262 are filled as one
# contiguous range so Excel groups them into a single shared formula,
# mirroring the pattern used by every earlier block in this sheet.
$ws.Range("I250").Formula = '="static constexpr u8 "&C250&" = "&B250&";"'
$ws.Range("I251:I262").Formula = '="static constexpr u8 "&C251&" = "&B251&";"'

# Column J - registerField(...) call text, same row-grouping convention.
$ws.Range("J250").Formula = '="registerField(Fields::"&C250&", """&C250&""", "&D250&", "&E250&", "&F250&");"'
$ws.Range("J251:J262").Formula = '="registerField(Fields::"&C251&", """&C251&""", "&D251&", "&E251&", "&F251&");"'

# ---------------------------------------------------------------------------
# Misc view-state bookkeeping to mirror the end-state captured by the diff:
# the sheet's used range grew, the user scrolled further down, and moved the
# active selection to below the newly added data.
# ---------------------------------------------------------------------------
$ws.Range("E263").Select()
$excel.ActiveWindow.ScrollRow = 235
